# Fruta / hortaliza, semanal
# Insert a new weekly record at row 49 (pushing the existing rows 49-107 down
# to 50-108) and populate the new row with this week's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(49).Insert()

$ws.Range("A49").Value = 10
$ws.Range("B49").Value = "Vega Modelo de Temuco"
$ws.Range("C49").Value = "La Araucanía"
$ws.Range("D49").Value = 44413
$ws.Range("E49").Value = 9
$ws.Range("F49").Value = "Fruta"
$ws.Range("G49").Value = 100102
$ws.Range("H49").Value = "Cítricos"
$ws.Range("I49").Value = 100102006
$ws.Range("J49").Value = "Pomelo"
$ws.Range("K49").Value = "Start Ruby"
$ws.Range("L49").Value = "Primera"
$ws.Range("M49").Value = 65
$ws.Range("N49").Value = 12000
$ws.Range("O49").Value = 12000
$ws.Range("P49").Value = 12000
$ws.Range("Q49").Value = "`$/bandeja 15 kilos granel"
$ws.Range("R49").Value = "Región de O'Higgins"
$ws.Range("S49").Value = 800
$ws.Range("T49").Value = 15
